# Apply "modified div and final_div parameters" edit to Observation workbook
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fill in missing E and G values for rows 4 and 5 ---
$ws.Range("E4").Value = 86.45
$ws.Range("G4").Value = 86.25

$ws.Range("E5").Value = 86.45
$ws.Range("G5").Value = 86.25

# --- Add new row 6: div_factor = 1 ---
$ws.Range("A6").Value = 3
$ws.Range("B6").Value = "div_factor"
$ws.Range("C6").Value = 1

# --- Add new row 7: final_div_factor = 25 ---
$ws.Range("A7").Value = 3
$ws.Range("B7").Value = "final_div_factor"
$ws.Range("C7").Value = 25

# Match formatting of the other "Change_ID" group rows (2 and 3), which use
# the yellow-highlighted style, by copying that formatting onto the new rows.
$ws.Range("A2:C3").Copy() | Out-Null
$ws.Range("A6:C7").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# --- Update selection to match the final cursor position in the diff ---
$ws.Range("D7").Select() | Out-Null
